# Trade #18 closed at 2026-02-17 23:54:22 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" rollups for the
# MarketMaking strategy after its 18th trade closed, and appends the
# new trade row (#18) to both the "All Trades" and "MarketMaking"
# trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet - workbook-level rollup
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.32   # Current Capital
$summary.Range("B4").Value = 0.32      # Total P&L $
$summary.Range("B5").Value = 0.36      # Total P&L %
$summary.Range("B6").Value = 18        # Total Trades
$summary.Range("B8").Value = 7         # Losing Trades
$summary.Range("B9").Value = 55.56     # Win Rate %

# ---------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.32     # Capital
$status.Range("D6").Value = 18         # Trades
$status.Range("E6").Value = 0.32       # P&L $
$status.Range("F6").Value = 0.32       # P&L %
$status.Range("G6").Value = 55.56      # Win Rate %

# ---------------------------------------------------------------
# 3) Append new trade row (#18) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------
$newRow = 19

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 18                     # A: Trade #

    # Column B holds a plain date-like text ("2026-02-17"); a leading
    # apostrophe forces Excel to keep it as literal text instead of
    # auto-converting it to a date serial number.
    $ws.Cells.Item($newRow, 2).Value = "'2026-02-17"          # B: Date

    $ws.Cells.Item($newRow, 3).Value = "23:54:15"             # C: Time
    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"         # D: Strategy
    $ws.Cells.Item($newRow, 5).Value = "DOWN"                 # E: Side
    $ws.Cells.Item($newRow, 6).Value = 0.02                   # F: Entry Price
    $ws.Cells.Item($newRow, 7).Value = 0.01                   # G: Exit Price
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"               # H: Status
    $ws.Cells.Item($newRow, 9).Value = -50                    # I: P&L %
    $ws.Cells.Item($newRow, 10).Value = -0.01                 # J: P&L $
    $ws.Cells.Item($newRow, 11).Value = 100.32                # K: Capital After
    $ws.Cells.Item($newRow, 12).Value = 0                     # L: Entry Slippage (bps)
    $ws.Cells.Item($newRow, 13).Value = 0                     # M: Exit Slippage (bps)
    $ws.Cells.Item($newRow, 14).Value = 0.6                   # N: Confidence
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item($newRow, 16).Value = "early_exit"          # P: Exit Reason
    $ws.Cells.Item($newRow, 17).Value = 0.14                  # Q: Duration (min)
}
